# Applies "Add support for #ignore part way through lines in the framework"
# to the Characteristics sheet: the 'Setup Weight' column (E) trades places
# with 'Components'/'Denominator' (F/G), moving to the end (G) of that trio,
# a couple of new data points are added (H3=10, G9=#ignore, H9=1, blank row
# 10), and the active-sheet/selection bookkeeping is updated to match.

$wb = $excel.ActiveWorkbook

$wsChar   = $wb.Worksheets.Item("Characteristics")
$wsParams = $wb.Worksheets.Item("Parameters")

# ---------------------------------------------------------------------
# 1. Header row: Setup Weight (was E1) moves to G1; Components (was F1)
#    moves to E1; Denominator (was G1) moves to F1. H1 (Default Value)
#    is untouched.
# ---------------------------------------------------------------------
$wsChar.Range("E1").Value = "Components"
$wsChar.Range("F1").Value = "Denominator"
$wsChar.Range("G1").Value = "Setup Weight"

# ---------------------------------------------------------------------
# 2. Data rows 2-9: shift the old Components values (col F) left into
#    E, and the old Denominator values (col G) left into F, then blank
#    out whatever's left dangling in F/G (ClearContents keeps the
#    cell's existing style).
#    Column E previously held the (centred, style "2") Setup Weight
#    numbers; it now holds the (left-aligned, style "3") Components
#    text, so its alignment needs to follow the data.
# ---------------------------------------------------------------------
$wsChar.Range("E2:E9").HorizontalAlignment = -4131

$componentsByRow = @{
    2 = "sus, inf, rec"
    3 = "inf"
    4 = "inf, rec"
    5 = "inf, sus"
    6 = "sus, rec"
    7 = "ch_infrec"
    8 = "ch_infsus"
    9 = "ch_newinf"
}
$denominatorByRow = @{
    3 = "ch_all"
    7 = "ch_all"
    8 = "ch_all"
    9 = "ch_all"
}

foreach ($row in 2..9) {
    $wsChar.Range("E$row").Value = $componentsByRow[$row]

    if ($denominatorByRow.ContainsKey($row)) {
        $wsChar.Range("F$row").Value = $denominatorByRow[$row]
    } else {
        $wsChar.Range("F$row").ClearContents()
    }

    $wsChar.Range("G$row").ClearContents()
}

# New Setup Weight values that land in column G.
$wsChar.Range("G9").Value = "#ignore"
$wsChar.Range("H3").Value = 10
$wsChar.Range("H9").Value = 1

# New (otherwise empty) row 10, carrying only a left-aligned blank G10
# cell, matching the style already used throughout column G.
$wsChar.Range("G10").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 3. Comments: re-point E1/F1/G1 comments at their new columns. Excel
#    re-sorts the comment list by cell reference on save, so we don't
#    need to worry about insertion order here.
# ---------------------------------------------------------------------
$componentsComment = @"
This column, and any that immediately follow without a specified
header, is for the 'components' of a cascade characteristic.
A component is either a compartment or a characteristic that has
been previously defined, i.e. in a previous row, and should be
listed in this (and appropriate subsequent columns) by 'Code Name'.
For example, characteristic 'infected' may include 'dis_stage_1',
'dis_stage_2' and 'dis_advanced', where 'dis_advanced' is another
previously-defined characteristic including 'dis_stage_3' and
'dis_stage_4'.
In an example model, 'infected' would track population size summed
across the four 'dis_stage' states.
Note: If two or more components are listed in the same column, they
must be separated by a comma.
Whitespace is allowable and will be deleted during processing.
[attribute_charac_includes]
"@

$denominatorComment = "This column defines a 'denominator' attribute for a 'charac' item."

$setupWeightComment = @"
This column determines how important user-provided values for this
characteristic are to setting up the initial state of a model.
In general, the column value should be '1' if model construction
directly depends on what the user provides for the characteristic.
It should be '0' if supplied values are only for calibration or
note-keeping purposes.
In this latter case, the linear-algebra method of setting up
compartment sizes may complain about an 'under-determined' system
during a model run.
To avoid this, every compartment that has a nonzero setup weight
should be included in at least one distinct characteristic.

Note: Default value, i.e. a blank cell, is '1'.
Framework file parsing should also warn the user about a characteristic
with nonzero setup weight that is suppressed in the databook, i.e.
has a databook order of '-1'.
"@

$oldSetupWeightText = $wsChar.Range("E1").Comment.Text()

$wsChar.Range("E1").Comment.Delete()
$wsChar.Range("F1").Comment.Delete()
$wsChar.Range("G1").Comment.Delete()

$wsChar.Range("E1").AddComment($componentsComment.TrimEnd("`r", "`n"))
$wsChar.Range("F1").AddComment($denominatorComment)
$wsChar.Range("G1").AddComment($setupWeightComment.TrimEnd("`r", "`n"))

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping: the Parameters sheet's cursor
#    moved to G9, and the active tab moved from Transitions to
#    Characteristics (whose cursor ends on L14). Order matters: we
#    touch Parameters first so that the final Activate()/Select() on
#    Characteristics is what sticks as the active tab.
# ---------------------------------------------------------------------
$wsParams.Range("G9").Select()

$wsChar.Activate()
$wsChar.Range("L14").Select()
